$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B18
$ws.Range("B18").Value = 26394

# Fill in the new rows 19-23 with data (columns A through I)
$data = @(
    @(18, 24734, 0, 16117, 15506, 0, 0, 0, 0),
    @(19, 24108, 0, 15648, 15435, 0, 0, 0, 0),
    @(20, 23135, 0, 14515, 14707, 0, 0, 0, 0),
    @(21, 22359, 0, 13506, 14749, 0, 0, 0, 0),
    @(22, 22099, 0, 11221, 14631, 0, 0, 0, 0)
)

$startRow = 19
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowValues[$col - 1]
    }
}
